$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ETH")
$ws.Range("J3").Value = 2902.62042703556
$ws.Range("B12").Value = 0.00733564
$ws.Range("B35").Value = 0.12679807
$ws.Range("D35").Value = 221.64
$ws.Range("B36").Value = 0.02545015
$ws.Range("D36").Value = 45.1
$ws.Range("B40").Value = 0.05763159
$ws.Range("D40").Value = 109.05
$ws = $wb.Worksheets.Item("APE")
$ws.Range("J3").Value = 1.675114566754455
$ws.Range("B5").Value = 16.68810811
$ws.Range("D5").Value = 45.1
$ws.Range("B6").Value = 0.60078517
$ws = $wb.Worksheets.Item("ATOM")
$ws.Range("J3").Value = 9.839757168182871
$ws.Range("B7").Value = 0.02935176
$ws = $wb.Worksheets.Item("AVAX")
$ws.Range("J3").Value = 36.86705442985877
$ws.Range("B5").Value = 2.67501715
$ws.Range("D5").Value = 45.1
$ws.Range("B6").Value = 0.01674849
$ws.Range("B10").Value = 0.75996743
$ws.Range("D10").Value = 12.42
$ws = $wb.Worksheets.Item("AMP")
$ws.Range("J3").Value = 0.004187903369712374
$ws = $wb.Worksheets.Item("BNB")
$ws.Range("J3").Value = 375.1725608647923
$ws.Range("B10").Value = 0.0027559
$ws.Range("B11").Value = 0.58452141
$ws.Range("D11").Value = 166.77
$ws.Range("B12").Value = 0.15663043
$ws.Range("D12").Value = 45.1
$ws = $wb.Worksheets.Item("DOGE")
$ws.Range("J3").Value = 0.08347138603017691
$ws.Range("B6").Value = 0.29236126
$ws = $wb.Worksheets.Item("DOT")
$ws.Range("J3").Value = 7.354507422417237
$ws.Range("B5").Value = 7.83923568
$ws.Range("D5").Value = 45.1
$ws.Range("B6").Value = 0.08013128999999999
$ws = $wb.Worksheets.Item("EGLD")
$ws.Range("J3").Value = 56.57180543717337
$ws.Range("B6").Value = 0.00300022
$ws = $wb.Worksheets.Item("GRT")
$ws.Range("J3").Value = 0.2366996604036511
$ws = $wb.Worksheets.Item("ICP")
$ws.Range("J3").Value = 12.84191611145351
$ws.Range("B5").Value = 2.52050326
$ws.Range("D5").Value = 12.96
$ws.Range("B6").Value = 0.00236469
$ws = $wb.Worksheets.Item("BTC")
$ws.Range("J3").Value = 50997.95341795307
$ws.Range("B6").Value = 0.00035618
$ws.Range("B23").Value = 0.00758882
$ws.Range("D23").Value = 197.85
$ws.Range("B24").Value = 0.0016683
$ws.Range("D24").Value = 45.1
$ws.Range("B34").Value = 0.00212721
$ws.Range("D34").Value = 63.95
$ws = $wb.Worksheets.Item("KAVA")
$ws.Range("J3").Value = 0.7610537062559728
$ws = $wb.Worksheets.Item("LDO")
$ws.Range("J3").Value = 2.985908831371233
$ws.Range("B5").Value = 6.96025417
$ws.Range("D5").Value = 16.02
$ws.Range("B6").Value = 0.02037352
$ws = $wb.Worksheets.Item("LINK")
$ws.Range("J3").Value = 18.26468387933613
$ws.Range("B5").Value = 1.47992576
$ws.Range("D5").Value = 11.7
$ws.Range("B6").Value = 0.00248987
$ws = $wb.Worksheets.Item("LTC")
$ws.Range("J3").Value = 68.25227881597635
$ws.Range("B6").Value = 0.00134812
$ws = $wb.Worksheets.Item("LUNA")
$ws.Range("J3").Value = 0.6647540677099818
$ws.Range("B6").Value = 0.0585125
$ws = $wb.Worksheets.Item("LUNC")
$ws.Range("J3").Value = 0.0001220978403538877
$ws.Range("B18").Value = 5039.86241575
$ws = $wb.Worksheets.Item("MATIC")
$ws.Range("J3").Value = 0.9232205394646995
$ws.Range("B6").Value = 0.32960107
$ws.Range("B7").Value = 49.56779993
$ws.Range("D7").Value = 45.1
$ws = $wb.Worksheets.Item("MEME")
$ws.Range("J3").Value = 0.02505018332921594
$ws.Range("B6").Value = 0.06835887
$ws = $wb.Worksheets.Item("MINA")
$ws.Range("J3").Value = 1.272951803311154
$ws.Range("B6").Value = 0.35160407
$ws = $wb.Worksheets.Item("NEAR")
$ws.Range("J3").Value = 3.143667639192321
$ws.Range("B6").Value = 24.19388975
$ws.Range("D6").Value = 45.1
$ws.Range("B7").Value = 0.10315562
$ws = $wb.Worksheets.Item("SEI")
$ws.Range("J3").Value = 0.8475409733494793
$ws.Range("B6").Value = 0.07643529
$ws = $wb.Worksheets.Item("SHIB")
$ws.Range("J3").Value = 0.000009423697043374744
$ws.Range("B6").Value = 277.73
$ws = $wb.Worksheets.Item("SHPING")
$ws.Range("J3").Value = 0.00494863800916052
$ws = $wb.Worksheets.Item("SOL")
$ws.Range("J3").Value = 102.0284760499668
$ws.Range("B16").Value = 6.12055976
$ws.Range("D16").Value = 130.14
$ws.Range("B17").Value = 0.06485841000000001
$ws.Range("B18").Value = 1.92440451
$ws.Range("D18").Value = 45.1
$ws = $wb.Worksheets.Item("TRX")
$ws.Range("J3").Value = 0.1394023411830511
$ws.Range("B6").Value = 0.26650207
$ws = $wb.Worksheets.Item("UNI")
$ws.Range("J3").Value = 7.20747519117467
$ws.Range("B5").Value = 2.75575411
$ws.Range("D5").Value = 15.9
$ws.Range("B6").Value = 0.00275603
$ws = $wb.Worksheets.Item("XRP")
$ws.Range("J3").Value = 0.5441979280554681
$ws.Range("B6").Value = 0.8696042899999999
$ws = $wb.Worksheets.Item("TIA")
$ws.Range("J3").Value = 16.97403346365317
$ws.Range("B6").Value = 0.00429857
$ws = $wb.Worksheets.Item("DYDX")
$ws.Range("J3").Value = 2.917593425079952
$ws.Range("B6").Value = 0.00087644
$ws = $wb.Worksheets.Item("POLIS")
$ws.Range("J3").Value = 0.3869838861045061
$ws = $wb.Worksheets.Item("ATLAS")
$ws.Range("J3").Value = 0.004507825771214326
$ws = $wb.Worksheets.Item("ACE")
$ws.Range("J3").Value = 9.666921769647226
$ws.Range("B6").Value = 0.0000259
$ws = $wb.Worksheets.Item("ADA")
$ws.Range("J3").Value = 0.5853408639352949
$ws.Range("B6").Value = 0.78540654
$ws.Range("B7").Value = 124.110001
$ws.Range("D7").Value = 45.1
$ws = $wb.Worksheets.Item("ALGO")
$ws.Range("J3").Value = 0.1858948249494281
$ws.Range("B6").Value = 0.58538583
